$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pruned / re-sequenced one-hot class labels for rows 2-49 (columns A:C)
$data = @(
    @(1,0,0),
    @(1,0,0),
    @(0,1,0),
    @(1,0,0),
    @(1,0,0),
    @(1,0,0),
    @(0,1,0),
    @(0,1,0),
    @(0,1,0),
    @(0,1,0),
    @(0,0,1),
    @(0,0,1),
    @(0,1,0),
    @(1,0,0),
    @(1,0,0),
    @(1,0,0),
    @(0,1,0),
    @(0,0,1),
    @(1,0,0),
    @(0,0,1),
    @(0,0,1),
    @(0,1,0),
    @(1,0,0),
    @(0,0,1),
    @(0,0,1),
    @(1,0,0),
    @(1,0,0),
    @(0,1,0),
    @(0,0,1),
    @(0,1,0),
    @(0,1,0),
    @(0,0,1),
    @(1,0,0),
    @(1,0,0),
    @(1,0,0),
    @(0,0,1),
    @(0,1,0),
    @(0,1,0),
    @(0,0,1),
    @(0,1,0),
    @(1,0,0),
    @(0,0,1),
    @(0,1,0),
    @(0,0,1),
    @(1,0,0),
    @(0,1,0),
    @(1,0,0),
    @(0,1,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
